$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "Cat(s)"
$ws.Range("B8").Value = "Dog(s)"
$ws.Range("B9").Value = "Washer"
$ws.Range("B10").Value = "Kitchen"
$ws.Range("B11").Value = "Shampoo"
$ws.Range("B12").Value = "Hangers"
$ws.Range("B13").Value = "Heating"
$ws.Range("B14").Value = "Hot Tub"
$ws.Range("B15").Value = "Doorman"
$ws.Range("B19").Value = "Essentials"
$ws.Range("B20").Value = "Hair Dryer"
$ws.Range("B33").Value = "Suitable for Events"
$ws.Range("B34").Value = "Family/Kid Friendly"
$ws.Range("B35").Value = "Lock on Bedroom Door"
$ws.Range("B36").Value = "Elevator in Building"
